$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Range("A13").Value = "DeezNutz"
$ws.Range("B13").Value = "TestJob1!"
